$d = $word.ActiveDocument

# --- Step 1: "Ashkan bröt ner ..." paragraph -----------------------------
# The paragraph used to end with the literal text followed by a _GoBack
# bookmark and a trailing " " run. The new version drops the bookmark and
# simply ends the same run with a trailing space.
$ashkanPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Ashkan bröt ner html")) {
        $ashkanPara = $p
        break
    }
}
$ar = $ashkanPara.Range
# Exclude the trailing paragraph mark from the replace range.
$body = $d.Range($ar.Start, $ar.End - 1)
$body.Text = "Ashkan bröt ner html sidorna i mindre deler för marcus som sedan lagt in dem i php backend, detta så att vi kan jobba med varje del för sig utan att behöva ändra i alla filer hela tiden. Kollat på design idéer "

# --- Step 2: insert the Dag5 / Dag6 block after "Gruppen diskuterar..." --
$gruppenPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Gruppen diskuterar tillsammans")) {
        $gruppenPara = $p
    }
}
$insertionPoint = $d.Range($gruppenPara.Range.End, $gruppenPara.Range.End)

$blockXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:pBdr><w:bottom w:val="single" w:sz="6" w:space="1" w:color="auto"/></w:pBdr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Da</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>g5</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Fredag</w:t></w:r></w:p>
<w:p><w:r><w:t>Rikard arbetar med bug fixa slidern.</w:t></w:r></w:p>
<w:p><w:r><w:t>..</w:t></w:r></w:p>
<w:p/>
<w:p><w:pPr><w:pBdr><w:bottom w:val="single" w:sz="6" w:space="1" w:color="auto"/></w:pBdr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Da</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>g6</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Måndag</w:t></w:r></w:p>
<w:p><w:r><w:t xml:space="preserve">Rikard </w:t></w:r><w:r><w:t>”</w:t></w:r><w:r><w:t>löser</w:t></w:r><w:r><w:t>”</w:t></w:r><w:r><w:t xml:space="preserve"> bugen med sliden via att börja om på en ny slider, hjälper C.H. med att implmentera sina ändringar in i projektet.</w:t></w:r></w:p>
<w:p><w:r><w:t>..</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$insertionPoint.InsertXML($blockXml)
